# Mise a jour du fichier Excel public
# - "Prix Spot": nouvelle colonne AJ pour le 19-jul
# - "Gaz" / "CO2": nouvelle ligne 33 pour le 2025-07-17
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": new column AJ = "19-jul" -----------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy AI1's formatting (bold/centered/bordered header style) onto AJ1
# before writing its text, so the new header cell matches the others.
$ws1.Range("AI1").Copy($ws1.Range("AJ1"))
$ws1.Range("AJ1").Value = "19-jul"

$ws1.Range("AJ2").Value  = 98.8
$ws1.Range("AJ3").Value  = 90.79000000000001
$ws1.Range("AJ4").Value  = 79.70999999999999
$ws1.Range("AJ5").Value  = 53.41
$ws1.Range("AJ6").Value  = 49.57
$ws1.Range("AJ7").Value  = 48.62
$ws1.Range("AJ8").Value  = 47.99
$ws1.Range("AJ9").Value  = 56.25
$ws1.Range("AJ10").Value = 56.62
$ws1.Range("AJ11").Value = 34.46
$ws1.Range("AJ12").Value = 10
$ws1.Range("AJ13").Value = 0.65
$ws1.Range("AJ14").Value = 0
$ws1.Range("AJ15").Value = 0
$ws1.Range("AJ16").Value = 0
$ws1.Range("AJ17").Value = 2.37
$ws1.Range("AJ18").Value = 18.82
$ws1.Range("AJ19").Value = 19.38
$ws1.Range("AJ20").Value = 30.16
$ws1.Range("AJ21").Value = 42.32
$ws1.Range("AJ22").Value = 72.92
$ws1.Range("AJ23").Value = 94.17
$ws1.Range("AJ24").Value = 106.8
$ws1.Range("AJ25").Value = 96.25

# --- Sheet "Gaz": new row 33 = 2025-07-17 / 33.6 ------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

# Build the date label as plain text via a scratch formula cell (well
# outside the used range) and paste its computed value into A33 so it
# lands as a literal string - like every other date in the column -
# instead of being auto-converted into a date serial number. The scratch
# cell is cleared again right away so no trace of it is left behind.
$ws2.Range("D1").Formula = "=TEXT(DATE(2025,7,17),""yyyy-mm-dd"")"
$ws2.Range("D1").Copy()
$ws2.Range("A33").PasteSpecial(-4163)
$ws2.Range("D1").ClearContents()
$ws2.Range("B33").Value = 33.6

# --- Sheet "CO2": new row 33 = 2025-07-17 / 69.81999999999999 -----------
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("D1").Formula = "=TEXT(DATE(2025,7,17),""yyyy-mm-dd"")"
$ws3.Range("D1").Copy()
$ws3.Range("A33").PasteSpecial(-4163)
$ws3.Range("D1").ClearContents()
$ws3.Range("B33").Value = 69.81999999999999
